$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Weekday OD")
$ws.Rows("46:46").Delete()
$ws.Columns("AS:AS").Delete()
